$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.769.16'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.793.63'
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.09'
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.35'
$ws.Range("E6").Value = '  -1.12%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.517'
$ws.Range("E8").Value = '  -0.40%  '

$ws.Range("E9").Value = '  -0.54%  '

$ws.Range("E10").Value = '  +0.88%  '

$ws.Range("E11").Value = '  +2.72%  '

$ws.Range("E12").Value = '  -2.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.82'
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.432.60'
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.790.01'
$ws.Range("E15").Value = '  -0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.797.56'
$ws.Range("E16").Value = '  +0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.32'
$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("E18").Value = '  +1.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.05'
$ws.Range("E19").Value = '  +0.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '462.66'
$ws.Range("E20").Value = '  +0.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.80'
$ws.Range("E21").Value = '  -2.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.699'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000147'
$ws.Range("E23").Value = '  -6.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.76'
$ws.Range("E24").Value = '  -0.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.06'
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E26").Value = '  -0.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("E28").Value = '  -0.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.942.87'
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.44'
$ws.Range("E30").Value = '  +2.61%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.65'
$ws.Range("E31").Value = '  -4.27%  '

$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.17'
$ws.Range("E33").Value = '  -1.87%  '

$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.04'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0993'
$ws.Range("E36").Value = '  -0.78%  '

$ws.Range("E37").Value = '  +0.26%  '

$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.990'
$ws.Range("E38").Value = '  -0.31%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.23'
$ws.Range("E39").Value = '  -4.81%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.77'
$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '44.53'
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("E44").Value = '  -0.95%  '

$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.22'
$ws.Range("E46").Value = '  +1.45%  '

$ws.Range("E47").Value = '  +8.65%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '400.37'
$ws.Range("E48").Value = '  +1.43%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.35'
$ws.Range("E49").Value = '  +0.69%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.86'
$ws.Range("E50").Value = '  +1.98%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '27.15'
$ws.Range("E51").Value = '  +1.62%  '
